# Insert a new data row before current row 383 (Femacal de La Calera - Ajo),
# shifting the existing rows 383:469 down to 384:470, then populate the
# newly-inserted row 383 with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 383; this pushes old row 383 (and below)
# down to row 384, and extends the used range to row 470.
$ws.Rows.Item(383).Insert()

# Populate the new row 383. Columns A,B,C,E,F,G,H,I,N,O,Q,R carry the same
# values as the record that is now on row 384 (i.e. unchanged from the
# original row 383); columns D,J,K,L,M,P take the new values.
$ws.Range("A383").Value = 3
$ws.Range("B383").Value = "Femacal de La Calera"
$ws.Range("C383").Value = "Coquimbo"
$ws.Range("D383").Value = 44754
$ws.Range("E383").Value = 5
$ws.Range("F383").Value = 100112003
$ws.Range("G383").Value = "Ajo"
$ws.Range("H383").Value = "Chino"
$ws.Range("I383").Value = "Primera"
$ws.Range("J383").Value = 85
$ws.Range("K383").Value = 20000
$ws.Range("L383").Value = 21000
$ws.Range("M383").Value = 20529
$ws.Range("N383").Value = '$/caja 10 kilos'
$ws.Range("O383").Value = "China"
$ws.Range("P383").Value = 2053
$ws.Range("Q383").Value = 10
$ws.Range("R383").Value = "Hortaliza"

# Match the date formatting used by the rest of column D.
$ws.Range("D383").NumberFormat = $ws.Range("D384").NumberFormat
